# close #187: Remove unnecessary column name in values and proportionality
#
# The sheet has an "id"/"nome" pair of header columns (A/B) followed by
# several blocks of proportionality values. The "nome" column (B), which
# held city names (Alta Floresta, Ariquemes, Cacoal, Cerejeiras), is no
# longer needed, so delete it entirely. Excel shifts every column to its
# right (C:M) one position to the left (B:L) automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("B").Delete()

# Match the author's resulting selection in the saved file.
[void]$ws.Range("F12").Select()
